$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New entry: issue #83 (the 7th generation cat-cat bundle), dated 11/7.
$ws.Range("A33").Value = "11/7"
$ws.Range("C33").Value = "第83期 第七代貓貓包"

# Match the author's final cursor position/selection after the edit.
$ws.Activate()
$ws.Range("C33").Select()
